$d = $word.ActiveDocument

# Phase 1: replace each unique original paragraph text with a unique placeholder token
$phase1 = @()
$phase1 += ,@("Tato rámcová dokumentace je uzavírána mezi Společností Legamedis Group a.s., [[ICO_1]], se sídlem V [[ADDRESS_1]], zastoupenou JUDr. [[PERSON_1]], (dále jen „Společnost“) a souborem jednotlivců, jejichž údaje jsou uvedeny níže a kteří jsou v dokumentu označováni jako „Subjekty“.", "@@SLOT0@@")
$phase1 += ,@("[[PERSON_2]] – „s [[PERSON_2]]“", "@@SLOT1@@")
$phase1 += ,@("[[PERSON_3]] – „o [[PERSON_3]]“", "@@SLOT2@@")
$phase1 += ,@("[[PERSON_4]] – „k [[PERSON_4]]“", "@@SLOT3@@")
$phase1 += ,@("[[PERSON_5]] – „pro [[PERSON_5]]“", "@@SLOT4@@")
$phase1 += ,@("[[PERSON_6]] – „s [[PERSON_6]]“", "@@SLOT5@@")
$phase1 += ,@("[[PERSON_7]] – „u [[PERSON_7]]“", "@@SLOT6@@")
$phase1 += ,@("[[PERSON_8]] – „od [[PERSON_8]]“", "@@SLOT7@@")
$phase1 += ,@("[[PERSON_9]] – „pro [[PERSON_9]]“", "@@SLOT8@@")
$phase1 += ,@("[[PERSON_10]] – „k [[PERSON_10]]“", "@@SLOT9@@")
$phase1 += ,@("[[PERSON_11]] – „s [[PERSON_11]]“", "@@SLOT10@@")
$phase1 += ,@("[[PERSON_12]] – „o [[PERSON_12]]“", "@@SLOT11@@")
$phase1 += ,@("[[PERSON_13]] – „k [[PERSON_13]]“", "@@SLOT12@@")
$phase1 += ,@("[[PERSON_14]] – „s [[PERSON_14]]“", "@@SLOT13@@")
$phase1 += ,@("[[PERSON_15]] – „u [[PERSON_15]]“", "@@SLOT14@@")
$phase1 += ,@("[[PERSON_16]] – „s [[PERSON_16]]“", "@@SLOT15@@")
$phase1 += ,@("[[PERSON_17]] – „o [[PERSON_17]]“", "@@SLOT16@@")
$phase1 += ,@("[[PERSON_18]] – „k [[PERSON_18]]“", "@@SLOT17@@")
$phase1 += ,@("[[PERSON_19]] – „u [[PERSON_19]]“", "@@SLOT18@@")
$phase1 += ,@("[[PERSON_20]] – „pro [[PERSON_20]]“", "@@SLOT19@@")
$phase1 += ,@("[[PERSON_21]] – „s [[PERSON_21]]“", "@@SLOT20@@")
$phase1 += ,@("[[PERSON_22]] – „k [[PERSON_22]]“", "@@SLOT21@@")
$phase1 += ,@("[[PERSON_23]] – „s [[PERSON_23]]“", "@@SLOT22@@")
$phase1 += ,@("[[PERSON_24]] – „o [[PERSON_24]]“", "@@SLOT23@@")
$phase1 += ,@("[[PERSON_25]] – „pro [[PERSON_25]]“", "@@SLOT24@@")
$phase1 += ,@("[[PERSON_26]] – „s [[PERSON_26]]“", "@@SLOT25@@")
$phase1 += ,@("[[PERSON_27]] – „k [[PERSON_27]]“", "@@SLOT26@@")
$phase1 += ,@("[[PERSON_28]] – „s [[PERSON_28]]“", "@@SLOT27@@")
$phase1 += ,@("[[PERSON_29]] – „o [[PERSON_29]]“", "@@SLOT28@@")
$phase1 += ,@("[[PERSON_30]] – „u [[PERSON_30]]“", "@@SLOT29@@")
$phase1 += ,@("[[PERSON_31]] – „k [[PERSON_31]]“", "@@SLOT30@@")
$phase1 += ,@("[[PERSON_32]] – „se [[PERSON_32]]“", "@@SLOT31@@")
$phase1 += ,@("[[PERSON_33]] – „u [[PERSON_33]]“", "@@SLOT32@@")
$phase1 += ,@("[[PERSON_34]] – „o [[PERSON_34]]“", "@@SLOT33@@")
$phase1 += ,@("[[PERSON_35]] – „s [[PERSON_35]]“", "@@SLOT34@@")
$phase1 += ,@("[[PERSON_36]] – „k [[PERSON_36]]“", "@@SLOT35@@")
$phase1 += ,@("[[PERSON_37]] – „od [[PERSON_37]]“", "@@SLOT36@@")
$phase1 += ,@("[[PERSON_38]] – „s [[PERSON_38]]“", "@@SLOT37@@")
$phase1 += ,@("[[PERSON_39]] – „u [[PERSON_39]]“", "@@SLOT38@@")
$phase1 += ,@("[[PERSON_40]] – „o [[PERSON_40]]“", "@@SLOT39@@")
$phase1 += ,@("[[PERSON_41]] – „k [[PERSON_41]]“", "@@SLOT40@@")
$phase1 += ,@("V těchto řízeních bylo jednáno např. s [[PERSON_2]], [[PERSON_5]], [[PERSON_24]] či [[PERSON_42]].", "@@SLOT41@@")
$phase1 += ,@("svědek [[PERSON_28]] (ve výpovědi označen jako „svědek Černého“),", "@@SLOT42@@")
$phase1 += ,@("poškozená [[PERSON_13]] („vyjádření [[PERSON_13]]“),", "@@SLOT43@@")
$phase1 += ,@("obžalovaný [[PERSON_4]] („obhajoba [[PERSON_4]]“),", "@@SLOT44@@")
$phase1 += ,@("znalkyně MUDr. [[PERSON_3]], soudní znalkyně v oboru psychiatrie,", "@@SLOT45@@")
$phase1 += ,@("právní zástupkyně JUDr. [[PERSON_31]], advokátka,", "@@SLOT46@@")
$phase1 += ,@("tlumočník [[PERSON_30]], zapsaný v seznamu tlumočníků.", "@@SLOT47@@")
$phase1 += ,@("Alergologické vyšetření č. ALG/2025/22751 provedené MUDr. [[PERSON_11]],", "@@SLOT48@@")
$phase1 += ,@("Neurologické testy č. NEU/2025/44119 provedené MUDr. [[PERSON_34]],", "@@SLOT49@@")
$phase1 += ,@("Oční vyšetření č. OFT/2023/11281 provedené MUDr. [[PERSON_29]].", "@@SLOT50@@")
$phase1 += ,@("Zvláštní pozornost byla věnována výsledkům [[PERSON_15]], [[PERSON_19]] a [[PERSON_40]].", "@@SLOT51@@")
$phase1 += ,@("mobil Samsung Galaxy S22, [[IMEI_1]],", "@@SLOT52@@")
$phase1 += ,@("Některé Subjekty poskytly technické přístupy pro řešení kauz:", "@@SLOT53@@")
$phase1 += ,@("právní cloud účet ID: LEX-ACC-88221 (spravovala [[PERSON_41]]),", "@@SLOT54@@")
$phase1 += ,@("[[PERSON_38]] („výslech [[PERSON_38]]“),", "@@SLOT55@@")
$phase1 += ,@("[[PERSON_35]] („výpověď [[PERSON_35]]“),", "@@SLOT56@@")
$phase1 += ,@("[[PERSON_33]] („záznam o výslechu [[PERSON_33]]“),", "@@SLOT57@@")
$phase1 += ,@("[[PERSON_18]] („výslech [[PERSON_18]]“).", "@@SLOT58@@")
$phase1 += ,@("PhDr. [[PERSON_29]] – psychologický posudek,", "@@SLOT59@@")
$phase1 += ,@("MUDr. [[PERSON_24]] – posudek z traumatologie,", "@@SLOT60@@")
$phase1 += ,@("Ing. [[PERSON_8]] – expertiza IT infrastruktury.", "@@SLOT61@@")
$phase1 += ,@("Tyto účty byly doloženy např. od [[PERSON_20]], [[PERSON_36]] nebo [[PERSON_26]].", "@@SLOT62@@")
$phase1 += ,@("[[PERSON_37]],", "@@SLOT63@@")
$phase1 += ,@("[[PERSON_42]],", "@@SLOT64@@")
$phase1 += ,@("[[PERSON_23]],", "@@SLOT65@@")
$phase1 += ,@("[[PERSON_10]].", "@@SLOT66@@")

foreach ($pair in $phase1) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "PHASE1 NOT FOUND: $old"
    }
}

# Phase 2: replace each placeholder token with the final new text
$phase2 = @()
$phase2 += ,@("@@SLOT0@@", "Tato rámcová dokumentace je uzavírána mezi [[PERSON_2]] Group a.s., [[ICO_1]], se sídlem V [[ADDRESS_1]], zastoupenou JUDr. [[PERSON_1]], (dále jen „Společnost“) a souborem jednotlivců, jejichž údaje jsou uvedeny níže a kteří jsou v dokumentu označováni jako „Subjekty“.")
$phase2 += ,@("@@SLOT1@@", "[[PERSON_3]] – „s [[PERSON_4]]“")
$phase2 += ,@("@@SLOT2@@", "[[PERSON_5]] – „o [[PERSON_6]]“")
$phase2 += ,@("@@SLOT3@@", "[[PERSON_7]] – „k [[PERSON_8]]“")
$phase2 += ,@("@@SLOT4@@", "[[PERSON_9]] – „pro [[PERSON_9]]“")
$phase2 += ,@("@@SLOT5@@", "[[PERSON_10]] – „s [[PERSON_10]]“")
$phase2 += ,@("@@SLOT6@@", "[[PERSON_11]] – „u [[PERSON_11]]“")
$phase2 += ,@("@@SLOT7@@", "[[PERSON_12]] – „od [[PERSON_12]]“")
$phase2 += ,@("@@SLOT8@@", "[[PERSON_13]] – „pro [[PERSON_13]]“")
$phase2 += ,@("@@SLOT9@@", "[[PERSON_14]] – „k [[PERSON_15]]“")
$phase2 += ,@("@@SLOT10@@", "[[PERSON_16]] – „s [[PERSON_16]]“")
$phase2 += ,@("@@SLOT11@@", "[[PERSON_17]] – „o [[PERSON_18]]“")
$phase2 += ,@("@@SLOT12@@", "[[PERSON_19]] – „k [[PERSON_20]]“")
$phase2 += ,@("@@SLOT13@@", "[[PERSON_21]] – „s [[PERSON_22]]“")
$phase2 += ,@("@@SLOT14@@", "[[PERSON_23]] – „u [[PERSON_23]]“")
$phase2 += ,@("@@SLOT15@@", "[[PERSON_24]] – „s [[PERSON_24]]“")
$phase2 += ,@("@@SLOT16@@", "[[PERSON_25]] – „o [[PERSON_25]]“")
$phase2 += ,@("@@SLOT17@@", "[[PERSON_26]] – „k [[PERSON_26]]“")
$phase2 += ,@("@@SLOT18@@", "[[PERSON_27]] – „u [[PERSON_27]]“")
$phase2 += ,@("@@SLOT19@@", "[[PERSON_28]] – „pro [[PERSON_29]]“")
$phase2 += ,@("@@SLOT20@@", "[[PERSON_30]] – „s [[PERSON_30]]“")
$phase2 += ,@("@@SLOT21@@", "[[PERSON_31]] – „k [[PERSON_31]]“")
$phase2 += ,@("@@SLOT22@@", "[[PERSON_32]] – „s [[PERSON_33]]“")
$phase2 += ,@("@@SLOT23@@", "[[PERSON_34]] – „o [[PERSON_35]]“")
$phase2 += ,@("@@SLOT24@@", "[[PERSON_36]] – „pro [[PERSON_36]]“")
$phase2 += ,@("@@SLOT25@@", "[[PERSON_37]] – „s [[PERSON_38]]“")
$phase2 += ,@("@@SLOT26@@", "[[PERSON_39]] – „k [[PERSON_40]]“")
$phase2 += ,@("@@SLOT27@@", "[[PERSON_41]] – „s [[PERSON_41]]“")
$phase2 += ,@("@@SLOT28@@", "[[PERSON_42]] – „o [[PERSON_43]]“")
$phase2 += ,@("@@SLOT29@@", "[[PERSON_44]] – „u [[PERSON_44]]“")
$phase2 += ,@("@@SLOT30@@", "[[PERSON_45]] – „k [[PERSON_45]]“")
$phase2 += ,@("@@SLOT31@@", "[[PERSON_46]] – „se [[PERSON_47]]“")
$phase2 += ,@("@@SLOT32@@", "[[PERSON_48]] – „u [[PERSON_48]]“")
$phase2 += ,@("@@SLOT33@@", "[[PERSON_49]] – „o [[PERSON_50]]“")
$phase2 += ,@("@@SLOT34@@", "[[PERSON_51]] – „s [[PERSON_51]]“")
$phase2 += ,@("@@SLOT35@@", "[[PERSON_52]] – „k [[PERSON_53]]“")
$phase2 += ,@("@@SLOT36@@", "[[PERSON_54]] – „od [[PERSON_54]]“")
$phase2 += ,@("@@SLOT37@@", "[[PERSON_55]] – „s [[PERSON_55]]“")
$phase2 += ,@("@@SLOT38@@", "[[PERSON_56]] – „u [[PERSON_56]]“")
$phase2 += ,@("@@SLOT39@@", "[[PERSON_57]] – „o [[PERSON_58]]“")
$phase2 += ,@("@@SLOT40@@", "[[PERSON_59]] – „k [[PERSON_59]]“")
$phase2 += ,@("@@SLOT41@@", "V těchto řízeních bylo jednáno např. s [[PERSON_3]], [[PERSON_9]], [[PERSON_34]] či [[PERSON_60]].")
$phase2 += ,@("@@SLOT42@@", "svědek [[PERSON_41]] (ve výpovědi označen jako „svědek Černého“),")
$phase2 += ,@("@@SLOT43@@", "poškozená [[PERSON_19]] („vyjádření [[PERSON_19]]“),")
$phase2 += ,@("@@SLOT44@@", "obžalovaný [[PERSON_7]] („obhajoba [[PERSON_7]]“),")
$phase2 += ,@("@@SLOT45@@", "znalkyně MUDr. [[PERSON_5]], soudní znalkyně v oboru psychiatrie,")
$phase2 += ,@("@@SLOT46@@", "právní zástupkyně JUDr. [[PERSON_45]], advokátka,")
$phase2 += ,@("@@SLOT47@@", "tlumočník [[PERSON_44]], zapsaný v seznamu tlumočníků.")
$phase2 += ,@("@@SLOT48@@", "Alergologické vyšetření č. ALG/2025/22751 provedené MUDr. [[PERSON_16]],")
$phase2 += ,@("@@SLOT49@@", "Neurologické testy č. NEU/2025/44119 provedené MUDr. [[PERSON_49]],")
$phase2 += ,@("@@SLOT50@@", "Oční vyšetření č. OFT/2023/11281 provedené MUDr. [[PERSON_42]].")
$phase2 += ,@("@@SLOT51@@", "Zvláštní pozornost byla věnována výsledkům [[PERSON_23]], [[PERSON_27]] a [[PERSON_57]].")
$phase2 += ,@("@@SLOT52@@", "mobil [[PERSON_61]] S22, [[IMEI_1]],")
$phase2 += ,@("@@SLOT53@@", "[[PERSON_62]] poskytly technické přístupy pro řešení kauz:")
$phase2 += ,@("@@SLOT54@@", "právní cloud účet ID: LEX-ACC-88221 (spravovala [[PERSON_59]]),")
$phase2 += ,@("@@SLOT55@@", "[[PERSON_55]] („výslech [[PERSON_55]]“),")
$phase2 += ,@("@@SLOT56@@", "[[PERSON_51]] („výpověď [[PERSON_51]]“),")
$phase2 += ,@("@@SLOT57@@", "[[PERSON_48]] („záznam o výslechu [[PERSON_48]]“),")
$phase2 += ,@("@@SLOT58@@", "[[PERSON_26]] („výslech [[PERSON_63]]“).")
$phase2 += ,@("@@SLOT59@@", "PhDr. [[PERSON_42]] – psychologický posudek,")
$phase2 += ,@("@@SLOT60@@", "MUDr. [[PERSON_34]] – posudek z traumatologie,")
$phase2 += ,@("@@SLOT61@@", "Ing. [[PERSON_12]] – expertiza IT infrastruktury.")
$phase2 += ,@("@@SLOT62@@", "Tyto účty byly doloženy např. od [[PERSON_28]], [[PERSON_52]] nebo [[PERSON_64]].")
$phase2 += ,@("@@SLOT63@@", "[[PERSON_54]],")
$phase2 += ,@("@@SLOT64@@", "[[PERSON_65]],")
$phase2 += ,@("@@SLOT65@@", "[[PERSON_32]],")
$phase2 += ,@("@@SLOT66@@", "[[PERSON_14]].")

foreach ($pair in $phase2) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "PHASE2 NOT FOUND: $old"
    }
}
Write-Output "DONE"
